$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.484.56"

$ws.Range("D3").Value = "1.619.45"
$ws.Range("E3").Value = "  -1.60%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'211.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.03%  "

$ws.Range("E6").Value = "  -1.36%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'22.84"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("E9").Value = "  +2.27%  "

$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").Value = "'0.0886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("D12").Value = "1.848.29"
$ws.Range("E12").Value = "  -1.63%  "

$ws.Range("D13").Value = "1.617.29"
$ws.Range("E13").Value = "  -1.71%  "

$ws.Range("E14").Value = "  -0.41%  "

$ws.Range("E15").Value = "  -2.56%  "

$ws.Range("D16").Value = "'65.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("D17").Value = "27.487.63"
$ws.Range("E17").Value = "  -0.60%  "

$ws.Range("D18").Value = "'230.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.23%  "

$ws.Range("D19").Value = "0.0₃0720"
$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("D20").Value = "'7.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.50%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("E22").Value = "  -0.77%  "

$ws.Range("D23").Value = "'10.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.25%  "

$ws.Range("E24").Value = "  +6.31%  "

$ws.Range("D25").Value = "'150.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.76%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.111"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.03%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'6.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.79%  "

$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("E30").Value = "  -0.87%  "

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("E32").Value = "  -1.05%  "

$ws.Range("D33").Value = "1.449.26"
$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("E34").Value = "  -3.48%  "

$ws.Range("E35").Value = "  -3.36%  "

$ws.Range("D36").Value = "'2.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("D37").Value = "'0.936"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.03%  "

$ws.Range("E38").Value = "  -1.89%  "

$ws.Range("E39").Value = "  -0.38%  "

$ws.Range("D40").Value = "'0.864"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.07%  "

$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("D42").Value = "'67.96"
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = "  +0.73%  "

$ws.Range("D44").Value = "'0.991"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.00%  "

$ws.Range("D45").Value = "'5.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.76%  "

$ws.Range("E46").Value = "  -2.21%  "

$ws.Range("D47").Value = "1.760.11"
$ws.Range("E47").Value = "  -1.55%  "

$ws.Range("E48").Value = "  +1.10%  "

$ws.Range("D49").Value = "'86.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("E50").Value = "  +18.27%  "

$ws.Range("D51").Value = "'0.101"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.72%  "
